$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.326636672019958
$ws.Range("B1").Value = 1.07689368724823
$ws.Range("C1").Value = 4.310553550720215
$ws.Range("D1").Value = 2.430655241012573
$ws.Range("E1").Value = 0.7423865795135498
